$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - numeric metrics (first set)
$ws.Range("B2").Value = 0.7911308369771046
$ws.Range("C2").Value = 0.8529849757961881
$ws.Range("D2").Value = 0.8843807794783712
$ws.Range("E2").Value = 0.6145278463798798
$ws.Range("F2").Value = 0.8336114287689569

# Row 3 - model descriptions (text)
$ws.Range("C3").Value = "RandomForestRegressor(bootstrap=False, max_features='log2', min_samples_split=8,`n                      n_estimators=20)"
$ws.Range("D3").Value = "XGBRegressor(base_score=0.5, booster='gbtree', colsample_bylevel=1,`n             colsample_bynode=1, colsample_bytree=1, gamma=0, gpu_id=-1,`n             importance_type='gain', interaction_constraints='',`n             learning_rate=0.1, max_delta_step=0, max_depth=2,`n             min_child_weight=3, missing=nan, monotone_constraints='()',`n             n_estimators=100, n_jobs=2, num_parallel_tree=1, random_state=0,`n             reg_alpha=0, reg_lambda=1, scale_pos_weight=1, subsample=1,`n             tree_method='exact', validate_parameters=1, verbosity=None)"
$ws.Range("E3").Value = "DecisionTreeRegressor(max_depth=6, min_samples_leaf=20, min_samples_split=20)"
$ws.Range("F3").Value = "AdaBoostRegressor(learning_rate=0.01, n_estimators=2000, random_state=1)"

# Row 4 - numeric metrics (second set)
$ws.Range("B4").Value = 51.28397524569497
$ws.Range("C4").Value = 36.0968309198763
$ws.Range("D4").Value = 28.38816969122341
$ws.Range("E4").Value = 94.64558625148524
$ws.Range("F4").Value = 40.85364849785898

# Row 5 - numeric metrics (third set)
$ws.Range("B5").Value = 5.509141880637524
$ws.Range("C5").Value = 4.626085546398045
$ws.Range("D5").Value = 3.450379840899736
$ws.Range("E5").Value = 7.147151009713077
$ws.Range("F5").Value = 4.950831660752229
